$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet tab/title to reflect the new "through" date
$ws.Name = "Through 2022-11-08"

# Update the "November (through 11-07)" label -> "November (through 11-08)"
$ws.Range("A12").Value = "November (through 11-08)"

# Update September 2022 value (I9)
$ws.Cells.Item(9, 9).Value = 162

# Update November row (row 12) values for years 2015-2022 (columns B-I)
$ws.Cells.Item(12, 2).Value = 10
$ws.Cells.Item(12, 3).Value = 19
$ws.Cells.Item(12, 4).Value = 29
$ws.Cells.Item(12, 5).Value = 20
$ws.Cells.Item(12, 6).Value = 12
$ws.Cells.Item(12, 7).Value = 50
$ws.Cells.Item(12, 8).Value = 57
$ws.Cells.Item(12, 9).Value = 22

# Update Total row (row 13) values for years 2015-2022 (columns B-I)
$ws.Cells.Item(13, 2).Value = 268
$ws.Cells.Item(13, 3).Value = 505
$ws.Cells.Item(13, 4).Value = 739
$ws.Cells.Item(13, 5).Value = 635
$ws.Cells.Item(13, 6).Value = 494
$ws.Cells.Item(13, 7).Value = 1107
$ws.Cells.Item(13, 8).Value = 1498
$ws.Cells.Item(13, 9).Value = 1421
